$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values
$ws.Range("B2").Value = 22.238109643292994
$ws.Range("C2").Value = -2.5972647720588284
$ws.Range("D2").Value = 0.39925577643444399
$ws.Range("E2").Value = 3.3322129417138484

# Row 3 values
$ws.Range("B3").Value = 19.973490344849282
$ws.Range("C3").Value = 3.6266888763321674
$ws.Range("D3").Value = -13.012829238525661
$ws.Range("E3").Value = 6.9050263606665823

# Update selection to match new sqref
$ws.Range("B1:E3").Select() | Out-Null
